$d = $word.ActiveDocument

$pairs = @(
    @{old="97×33=3201"; new="39×85=3315"},
    @{old="34×95=3230"; new="66×92=6072"},
    @{old="67×13=871";  new="45×93=4185"},
    @{old="66×31=2046"; new="81×40=3240"},
    @{old="72×80=5760"; new="64×13=832"},
    @{old="25×37=925";  new="64×39=2496"},
    @{old="93×60=5580"; new="30×42=1260"},
    @{old="49×26=1274"; new="93×26=2418"},
    @{old="47×28=1316"; new="75×64=4800"},
    @{old="84×32=2688"; new="42×61=2562"},
    @{old="39×16=624";  new="81×31=2511"},
    @{old="85×78=6630"; new="22×25=550"},
    @{old="63×51=3213"; new="89×46=4094"},
    @{old="21×36=756";  new="45×30=1350"},
    @{old="21×98=2058"; new="31×74=2294"},
    @{old="98×54=5292"; new="12×58=696"},
    @{old="50×36=1800"; new="96×62=5952"},
    @{old="47×12=564";  new="60×46=2760"},
    @{old="82×32=2624"; new="74×28=2072"},
    @{old="30×87=2610"; new="35×25=875"},
    @{old="44×49=2156"; new="69×14=966"},
    @{old="58×20=1160"; new="22×69=1518"},
    @{old="61×20=1220"; new="63×65=4095"},
    @{old="26×82=2132"; new="35×87=3045"},
    @{old="66×24=1584"; new="56×76=4256"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
